# Add partial project-data parse results into column B of the Template
# sheet. The rows correspond to PROJECT fields already labeled in column A
# (submitter id, availability type, code, investigator affiliation,
# investigator name, name, short name, support source, support id, state);
# "date collected" (row 20) is left blank since it wasn't parsed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B17").Value = "proj_sub_id"
$ws.Range("B18").Value = "Restricted"
$ws.Range("B19").Value = "code"
$ws.Range("B21").Value = "affil"
$ws.Range("B22").Value = "inv name"
$ws.Range("B23").Value = "name"
$ws.Range("B24").Value = "short name"
$ws.Range("B25").Value = "support source"
$ws.Range("B26").Value = "support id"
$ws.Range("B27").Value = "open"

# B26's row has no row-level format, so the new value wouldn't otherwise
# pick up the "field value" look used by the rest of the PROJECT block
# (rows 17-25). Paint that formatting on explicitly, matching a
# fill-down/format-paint from the row above.
$ws.Range("B24").Copy()
$ws.Range("B26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the selection on the last-edited cell, matching the author's
# workflow of filling the PROJECT section down to the state field.
[void]$ws.Range("B27").Select()
